$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.554.21"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.584.82"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'208.53"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'22.44"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").Value = "'0.0591"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "'0.0866"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "1.809.13"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "1.571.38"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "'0.526"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "27.573.61"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'63.20"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "'215.91"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.34"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0692"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "'153.46"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "'6.95"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'15.08"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").Value = "'1.16"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "'0.0473"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").Value = "1.372.08"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'2.96"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'1.55"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "'0.974"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("D39").Value = "'0.534"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").Value = "'0.828"
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'0.971"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.24"
$ws.Range("E43").Value = "  +6.31%  "
$ws.Range("D44").Value = "'1.80"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'64.55"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'5.30"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "1.720.68"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "'85.55"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'0.0961"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").Value = "'0.0495"
$ws.Range("E51").Value = "  -0.64%  "
